$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared string change: stacked_lstm -> bidirectional_lstm
$ws.Range("W3").Value = "bidirectional_lstm"

# Row 2
$ws.Range("B2").Value = 0.006527
$ws.Range("C2").Value = 0.004222
$ws.Range("D2").Value = 75
$ws.Range("E2").Value = -1.535421
$ws.Range("F2").Value = 0.12468
$ws.Range("G2").Value = 9.456235
$ws.Range("H2").Value = 54.887629
$ws.Range("I2").Value = -0.003357
$ws.Range("J2").Value = 9.094822000000001
$ws.Range("K2").Value = 75
$ws.Range("L2").Value = 3.031607
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 2

# Row 3
$ws.Range("B3").Value = 0.007225
$ws.Range("C3").Value = 0.004847
$ws.Range("D3").Value = 66.666667
$ws.Range("E3").Value = -0.776621
$ws.Range("F3").Value = 0.437383
$ws.Range("G3").Value = 5.143012
$ws.Range("H3").Value = 15.873702
$ws.Range("I3").Value = -0.009551
$ws.Range("J3").Value = 2.895164
$ws.Range("K3").Value = 66.666667
$ws.Range("L3").Value = 1.447582
$ws.Range("X3").Value = 32
$ws.Range("Z3").Value = 10
$ws.Range("AA3").Value = 50
$ws.Range("AD3").Value = 96

# Row 4
$ws.Range("B4").Value = 0.008109
$ws.Range("C4").Value = 0.005323
$ws.Range("D4").Value = 66.666667
$ws.Range("E4").Value = 0.6175349999999999
$ws.Range("F4").Value = 0.536882
$ws.Range("G4").Value = -0.163401
$ws.Range("H4").Value = -0.126396
$ws.Range("I4").Value = -0.024288
$ws.Range("J4").Value = 0.967888
$ws.Range("K4").Value = 66.666667
$ws.Range("L4").Value = 0.483944
$ws.Range("S4").Value = 1

# Row 5
$ws.Range("B5").Value = 0.007679
$ws.Range("C5").Value = 0.005059
